$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "FBPixels" (sheet1.xml): add two new test rows before the
# trailing "End" marker row, and clear the row-highlight formatting
# that Excel stops re-emitting once these rows are touched by hand.
# ------------------------------------------------------------------
$wsFb = $wb.Worksheets.Item("FBPixels")

# Insert two fresh rows right above the "End" row (currently row 16)
# while row 15 still carries its original style, so the new rows
# inherit the A:C highlight style the diff keeps for them.
$wsFb.Rows.Item(16).Insert()
$wsFb.Rows.Item(16).Insert()

$wsFb.Range("A16").Value = "QA"
$wsFb.Range("B16").Value = "Sub-D"
$wsFb.Range("C16").Value = "deluxe25offp"
$wsFb.Range("D16").Value = "CCFlow"
$wsFb.Range("E16").Value = "Facebook"

$wsFb.Range("A17").Value = "QA"
$wsFb.Range("B17").Value = "Sub-D"
$wsFb.Range("C17").Value = "cpcb2017"
$wsFb.Range("D17").Value = "CCFlow"
$wsFb.Range("E17").Value = "Facebook"

# Rows 4-15 lose their fill/style entirely in the new file.
$wsFb.Range("A4:E15").ClearFormats()

# The new rows keep the A:C highlight style but D:E come in unstyled.
$wsFb.Range("D16:E17").ClearFormats()

# The "End" marker (now row 18) also loses its style.
$wsFb.Range("A18").ClearFormats()

# ------------------------------------------------------------------
# Sheet "FB Pixel Test data" (sheet4.xml): the author just scrolled /
# re-selected a different range while reviewing the buyflow output;
# move the live selection there without leaving it the active tab.
# ------------------------------------------------------------------
$wsFbTest = $wb.Worksheets.Item("FB Pixel Test data")
[void]$wsFbTest.Range("A117:C117").Select()

# Restore FBPixels as the active/selected tab (matches the saved file).
$wsFb.Activate()
